$d = $word.ActiveDocument

# --- Locate the target run -------------------------------------------------
# The document opens with a title paragraph made of two runs:
#   run 1: "Godot "              (note the trailing space)
#   run 2: "Game Requirements"
# Both runs share the same character formatting (Times New Roman, 28 half-
# point size run font / 14pt). The edit re-types the first run as three
# separate runs that together still read "Godot ":
#   " "  +  "Godot"  +  " "
# all three keeping the exact same formatting as before.

$oldRange = $d.Range(0, 6)       # "Godot " -- the very start of the document
if ($oldRange.Text -ne "Godot ") {
    throw "Unexpected document start: [$($oldRange.Text)]"
}

# --- Helper: apply the run's original character formatting -----------------
function Set-TitleRunFont($range) {
    $range.Font.NameAscii = "Times New Roman"
    $range.Font.NameFarEast = "Times New Roman"
    $range.Font.NameOther = "Times New Roman"
    $range.Font.NameBi = "Times New Roman"
    $range.Font.Size = 14
    $range.Font.SizeBi = 14
}

# --- Perform the retype under tracked changes -------------------------------
# Plain (untracked) Range.InsertBefore/Delete calls get silently re-merged
# into a single run by the engine whenever two adjacent runs end up with
# identical formatting, which would collapse our three new runs back into
# one. Wrapping the edit in tracked changes (and accepting immediately
# after) keeps the three runs distinct in the saved XML, matching a real
# "select word, retype it" edit made with Track Changes off but still
# going through Word's insert/delete machinery.
$origTrack = $d.TrackRevisions
$d.TrackRevisions = $true

$oldRange.Delete()

# Insert the three replacement runs. Because the just-deleted text is kept
# around (hidden) as a tracked deletion, new text inserted in the "middle"
# of that hidden span gets relocated to its far edge. To keep everything
# landing in the right visible order we always insert at the very start of
# the document (position 0) and build the replacement back-to-front.
$insTrailingSpace = $d.Range(0, 0)
$insTrailingSpace.InsertBefore(" ")
Set-TitleRunFont($d.Range(0, 1))

$insWord = $d.Range(0, 0)
$insWord.InsertBefore("Godot")
Set-TitleRunFont($d.Range(0, 5))

$insLeadingSpace = $d.Range(0, 0)
$insLeadingSpace.InsertBefore(" ")
Set-TitleRunFont($d.Range(0, 1))

$d.TrackRevisions = $false
$d.AcceptAllRevisions()
$d.TrackRevisions = $origTrack

# Sanity check: visible text must be unchanged.
$check = $d.Range(0, 7)
if ($check.Text -ne " Godot ") {
    throw "Post-edit text mismatch: [$($check.Text)]"
}
